$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 252.77777
$ws.Range("I33").Value = 221.75
$ws.Range("K33").Value = 221.75
$ws.Range("M33").Value = 7.25

$ws.Range("H98").Value = 7431.56
$ws.Range("I98").Value = 4310.8887
$ws.Range("J98").Value = 9186.9375
$ws.Range("K98").Value = 4310.8887
$ws.Range("L98").Value = 9186.9375
$ws.Range("M98").Value = -2812.8887
$ws.Range("N98").Value = -12182.9375

$ws.Range("H108").Value = 39800
$ws.Range("J108").Value = 39800
$ws.Range("L108").Value = 39800
$ws.Range("N108").Value = -47480

$ws.Range("H109").Value = 38500
$ws.Range("J109").Value = 38500
$ws.Range("L109").Value = 38500
$ws.Range("N109").Value = -41274

$ws.Range("H112").Value = 1327.2642
$ws.Range("J112").Value = 1327.2642
$ws.Range("L112").Value = 3981.7926
$ws.Range("N112").Value = -6197.792600000001

$ws.Range("H113").Value = 5561.875
$ws.Range("I113").Value = 3015.8333
$ws.Range("J113").Value = 13200
$ws.Range("K113").Value = 3015.8333
$ws.Range("L113").Value = 13200
$ws.Range("M113").Value = 238.1667000000002
$ws.Range("N113").Value = -19708

$ws.Range("H122").Value = 7431.56
$ws.Range("I122").Value = 4310.8887
$ws.Range("J122").Value = 9186.9375
$ws.Range("K122").Value = 12932.6661
$ws.Range("L122").Value = 27560.8125
$ws.Range("M122").Value = -10482.6661
$ws.Range("N122").Value = -32460.8125

$ws.Range("H129").Value = 844.8
$ws.Range("J129").Value = 865.44794
$ws.Range("L129").Value = 2596.34382
$ws.Range("N129").Value = -12596.34382

$ws.Range("H135").Value = 719.3333
$ws.Range("I135").Value = 448.2
$ws.Range("J135").Value = 2075
$ws.Range("K135").Value = 4033.8
$ws.Range("L135").Value = 18675
$ws.Range("M135").Value = -1498.8
$ws.Range("N135").Value = -23745

$ws.Range("H141").Value = 57754.668
$ws.Range("I141").Value = 68486.266
$ws.Range("K141").Value = 205458.798
$ws.Range("M141").Value = -200278.798

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3501.25
$ws.Range("I45").Value = 4003.6667
$ws.Range("K45").Value = 4003.6667
$ws.Range("M45").Value = -3626.6667

$ws.Range("H61").Value = 1729.625
$ws.Range("I61").Value = 1651.1666
$ws.Range("K61").Value = 1651.1666
$ws.Range("M61").Value = -1439.1666

$ws.Range("H74").Value = 1282.2727
$ws.Range("I74").Value = 621.0526
$ws.Range("J74").Value = 5470
$ws.Range("K74").Value = 621.0526
$ws.Range("L74").Value = 5470
$ws.Range("M74").Value = 252.9474
$ws.Range("N74").Value = -7218

$ws.Range("H77").Value = 1282.2727
$ws.Range("I77").Value = 621.0526
$ws.Range("J77").Value = 5470
$ws.Range("K77").Value = 3105.263
$ws.Range("L77").Value = 27350
$ws.Range("M77").Value = 1262.737
$ws.Range("N77").Value = -36086

$ws.Range("H110").Value = 837.8
$ws.Range("I110").Value = 847.25
$ws.Range("J110").Value = 800
$ws.Range("K110").Value = 847.25
$ws.Range("L110").Value = 800
$ws.Range("M110").Value = 1197.75
$ws.Range("N110").Value = -4890

$ws.Range("H132").Value = 2248.28
$ws.Range("I132").Value = 1105.0625
$ws.Range("K132").Value = 3315.1875
$ws.Range("M132").Value = -785.1875

$ws.Range("H136").Value = 1729.625
$ws.Range("I136").Value = 1651.1666
$ws.Range("K136").Value = 4953.4998
$ws.Range("M136").Value = -2403.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2185.5625
$ws.Range("I107").Value = 1966.9
$ws.Range("J107").Value = 2550
$ws.Range("K107").Value = 1966.9
$ws.Range("L107").Value = 2550
$ws.Range("M107").Value = -46.90000000000009
$ws.Range("N107").Value = -6390

$ws.Range("H134").Value = 3429
$ws.Range("I134").Value = 1741.5333
$ws.Range("J134").Value = 9757
$ws.Range("K134").Value = 5224.5999
$ws.Range("L134").Value = 29271
$ws.Range("M134").Value = -2689.5999
$ws.Range("N134").Value = -34341

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6960.593
$ws.Range("I31").Value = 1423.6364
$ws.Range("J31").Value = 10767.25
$ws.Range("K31").Value = 1423.6364
$ws.Range("L31").Value = 10767.25
$ws.Range("M31").Value = -1128.6364
$ws.Range("N31").Value = -11357.25

$ws.Range("H34").Value = 6960.593
$ws.Range("I34").Value = 1423.6364
$ws.Range("J34").Value = 10767.25
$ws.Range("K34").Value = 1423.6364
$ws.Range("L34").Value = 10767.25
$ws.Range("M34").Value = -1221.6364
$ws.Range("N34").Value = -11171.25

$ws.Range("H58").Value = 1985.5938
$ws.Range("I58").Value = 1544.25
$ws.Range("J58").Value = 5075
$ws.Range("K58").Value = 1544.25
$ws.Range("L58").Value = 5075
$ws.Range("M58").Value = -1341.25
$ws.Range("N58").Value = -5481

$ws.Range("H136").Value = 1985.5938
$ws.Range("I136").Value = 1544.25
$ws.Range("J136").Value = 5075
$ws.Range("K136").Value = 4632.75
$ws.Range("L136").Value = 15225
$ws.Range("M136").Value = -2082.75
$ws.Range("N136").Value = -20325

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3111.0952
$ws.Range("I122").Value = 870.8889
$ws.Range("J122").Value = 3722.0605
$ws.Range("K122").Value = 7838.0001
$ws.Range("L122").Value = 33498.5445
$ws.Range("M122").Value = -5388.0001
$ws.Range("N122").Value = -38398.5445

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 181.875
$ws.Range("I2").Value = 55.454544
$ws.Range("J2").Value = 460
$ws.Range("K2").Value = 55.454544
$ws.Range("L2").Value = 460
$ws.Range("M2").Value = 57.545456
$ws.Range("N2").Value = -686

$ws.Range("H102").Value = 2480
$ws.Range("I102").Value = 1241.8182
$ws.Range("K102").Value = 1241.8182
$ws.Range("M102").Value = 380.1818000000001

$ws.Range("H113").Value = 1820
$ws.Range("I113").Value = 1525
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1525
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 645
$ws.Range("N113").Value = -7340

$ws.Range("H123").Value = 10325.934
$ws.Range("J123").Value = 10325.934
$ws.Range("L123").Value = 10325.934
$ws.Range("N123").Value = -15225.934

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1103.8306
$ws.Range("I68").Value = 966.2909
$ws.Range("J68").Value = 2995
$ws.Range("K68").Value = 966.2909
$ws.Range("L68").Value = 2995
$ws.Range("M68").Value = -217.2909
$ws.Range("N68").Value = -4493

$ws.Range("H71").Value = 1103.8306
$ws.Range("I71").Value = 966.2909
$ws.Range("J71").Value = 2995
$ws.Range("K71").Value = 4831.4545
$ws.Range("L71").Value = 14975
$ws.Range("M71").Value = -1087.4545
$ws.Range("N71").Value = -22463

$ws.Range("H133").Value = 20900.4
$ws.Range("J133").Value = 24456.615
$ws.Range("L133").Value = 24456.615
$ws.Range("N133").Value = -29516.615

$ws.Range("H136").Value = 3474.1052
$ws.Range("I136").Value = 1475.6666
$ws.Range("K136").Value = 4426.9998
$ws.Range("M136").Value = -1876.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 26111.111
$ws.Range("J64").Value = 26111.111
$ws.Range("L64").Value = 26111.111
$ws.Range("N64").Value = -26607.111

$ws.Range("H67").Value = 26111.111
$ws.Range("J67").Value = 26111.111
$ws.Range("L67").Value = 26111.111
$ws.Range("N67").Value = -27827.111
